$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string header "Corea.del.Sur" in G1
$ws.Cells.Item(1, 7).Value = "Corea.del.Sur"

# Fill in row 38 (23 de Marzo) which was previously only A38/B38
$ws.Cells.Item(38, 3).Value = 746
$ws.Cells.Item(38, 4).Value = 1546
$ws.Cells.Item(38, 5).Value = 28768
$ws.Cells.Item(38, 6).Value = 59138

# Populate new column G (Corea del Sur) for rows 2-38
$coreaValues = @{
    2 = 28
    3 = 29
    4 = 30
    5 = 31
    6 = 58
    7 = 111
    8 = 209
    9 = 436
    10 = 602
    11 = 833
    12 = 977
    13 = 1261
    14 = 1766
    15 = 2337
    16 = 3150
    17 = 3736
    18 = 4335
    19 = 5186
    20 = 5621
    21 = 6284
    22 = 6593
    23 = 7041
    24 = 7313
    25 = 7478
    26 = 7513
    27 = 7755
    28 = 7869
    29 = 7979
    30 = 8086
    31 = 8162
    32 = 8236
    33 = 8320
    34 = 8413
    35 = 8565
    36 = 8652
    37 = 8799
    38 = 8897
}
foreach ($r in $coreaValues.Keys) {
    $ws.Cells.Item($r, 7).Value = $coreaValues[$r]
}

# Update selection to match the author's final cursor position (G2)
$ws.Range("G2").Select()

Write-Output "done"
